$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$hls = $ws.Hyperlinks
$n = $hls.Count()
for ($i = 1; $i -le 5; $i++) {
    $h = $hls.Item($i)
    $rng = $h.Range()
    $addr = $rng.Address()
    $val = $rng.Value()
    Write-Host $i $addr "=>" $val
}
